$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Tipo" shifts from D -> E), carrying the
# header style along so the new D1 cell matches the existing header format.
$ws.Range("D1").EntireColumn.Insert()

# New "MAE" header in the freshly inserted column D.
$ws.Range("D1").Value = "MAE"

# Updated metric values in row 2.
$ws.Range("B2").Value = 0.3374190203127619
$ws.Range("C2").Value = 0.9935088083481195
$ws.Range("D2").Value = 0.4643319800914053
